# "After all fixes and flow"
#
# Rebuilds the workbook from the original 2-sheet "Opencart_LoginData" file
# into the final 4-sheet layout:
#   1. reg     (new)       - registration test data
#   2. login   (was Sheet1)- login test data (updated username/password)
#   3. address (new)       - address test data
#   4. query   (was Sheet2)- "contact us" query test data (unchanged)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Rename the two existing sheets in place, then append two new ones.
# ---------------------------------------------------------------------
$wsReg   = $wb.Worksheets.Item(1)
$wsReg.Name = "reg"

$wsLogin = $wb.Worksheets.Item(2)
$wsLogin.Name = "login"

$wsAddress = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsAddress.Name = "address"

$wsQuery = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsQuery.Name = "query"

# ---------------------------------------------------------------------
# 2) "reg" sheet - brand new registration data.
#    The original "Sheet1" had 5 data rows (username/password/res); the
#    new "reg" layout only needs 3 rows, so drop the trailing two.
# ---------------------------------------------------------------------
$wsReg.Range("A4:C5").EntireRow.Delete()

$wsReg.Cells.Item(1,1).Value = "Fname"
$wsReg.Cells.Item(1,2).Value = "Lname"
$wsReg.Cells.Item(1,3).Value = "email"
$wsReg.Cells.Item(1,4).Value = "telephone"
$wsReg.Cells.Item(1,5).Value = "pass"
$wsReg.Cells.Item(1,6).Value = "confirmpass"
$wsReg.Range("A1:F1").Font.Bold = $true

$wsReg.Cells.Item(2,1).Value = "Love"
$wsReg.Cells.Item(2,2).Value = "Patel"
$wsReg.Cells.Item(2,3).Value = "lovepatel1@yopmail.com"
$wsReg.Cells.Item(2,4).Value = 9411481766
$wsReg.Cells.Item(2,5).Value = "Nikk@12345"
$wsReg.Cells.Item(2,6).Value = "Nikk@12345"

$wsReg.Cells.Item(3,1).Value = "Himanshi"
$wsReg.Cells.Item(3,2).Value = "Patel"
$wsReg.Cells.Item(3,3).Value = "himanshi1@yopmail.com"
$wsReg.Cells.Item(3,4).Value = 9968198833
$wsReg.Cells.Item(3,5).Value = "Nikk@12345"
$wsReg.Cells.Item(3,6).Value = "Nikk@12345"

$wsReg.Columns.Item(3).ColumnWidth = 28.38
$wsReg.Columns.Item(6).ColumnWidth = 15.6

$wsReg.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&Kffffff&A'
$wsReg.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12&KffffffPage &P'
$wsReg.PageSetup.LeftMargin = 56.7
$wsReg.PageSetup.RightMargin = 56.7
$wsReg.PageSetup.TopMargin = 75.8
$wsReg.PageSetup.BottomMargin = 75.8
$wsReg.PageSetup.HeaderMargin = 56.7
$wsReg.PageSetup.FooterMargin = 56.7

# ---------------------------------------------------------------------
# 3) "login" sheet - replace the stale username/password rows.
# ---------------------------------------------------------------------
$wsLogin.Cells.Item(2,1).Value = "lovepatel1@yopmail.com"
$wsLogin.Cells.Item(2,2).Value = "Nikk@12345"
$wsLogin.Cells.Item(2,3).Value = "Valid"

$wsLogin.Cells.Item(4,1).Value = "pavan@gmail.com"
$wsLogin.Cells.Item(4,2).Value = "test123"
$wsLogin.Cells.Item(4,3).Value = "Invalid"

$wsLogin.Cells.Item(5,1).Value = "pavan@gmail.com"
$wsLogin.Cells.Item(5,2).Value = "test"
$wsLogin.Cells.Item(5,3).Value = "Invalid"

$wsLogin.Columns.Item(2).ColumnWidth = 22.67

# ---------------------------------------------------------------------
# 4) "address" sheet - brand new address-book data.
# ---------------------------------------------------------------------
$wsAddress.Cells.Item(1,1).Value = "Fname"
$wsAddress.Cells.Item(1,2).Value = "Lname"
$wsAddress.Cells.Item(1,3).Value = "company"
$wsAddress.Cells.Item(1,4).Value = "address1"
$wsAddress.Cells.Item(1,5).Value = "address2"
$wsAddress.Cells.Item(1,6).Value = "City"
$wsAddress.Cells.Item(1,7).Value = "Postcode"
$wsAddress.Range("A1:G1").Font.Bold = $true

$wsAddress.Cells.Item(2,1).Value = "Himanshi"
$wsAddress.Cells.Item(2,2).Value = "Patel"
$wsAddress.Cells.Item(2,3).Value = "Appinventiv"
$wsAddress.Cells.Item(2,4).Value = "bareilly"
$wsAddress.Cells.Item(2,5).Value = "delhi"
$wsAddress.Cells.Item(2,6).Value = "bareilly"
$wsAddress.Cells.Item(2,7).Value = 242503

$wsAddress.Cells.Item(3,1).Value = "Love"
$wsAddress.Cells.Item(3,2).Value = "Patel"
$wsAddress.Cells.Item(3,3).Value = "UPPP"
$wsAddress.Cells.Item(3,4).Value = "Noida"
$wsAddress.Cells.Item(3,5).Value = "Noida"
$wsAddress.Cells.Item(3,6).Value = "bareilly"
$wsAddress.Cells.Item(3,7).Value = 243507

$wsAddress.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&Kffffff&A'
$wsAddress.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12&KffffffPage &P'
$wsAddress.PageSetup.LeftMargin = 56.7
$wsAddress.PageSetup.RightMargin = 56.7
$wsAddress.PageSetup.TopMargin = 75.8
$wsAddress.PageSetup.BottomMargin = 75.8
$wsAddress.PageSetup.HeaderMargin = 56.7
$wsAddress.PageSetup.FooterMargin = 56.7

# ---------------------------------------------------------------------
# 5) "query" sheet (formerly Sheet2) - data is unchanged from before.
# ---------------------------------------------------------------------
$wsQuery.Cells.Item(1,1).Value = "name"
$wsQuery.Cells.Item(1,2).Value = "email"
$wsQuery.Cells.Item(1,3).Value = "query"

$wsQuery.Cells.Item(2,1).Value = "Niteesh"
$wsQuery.Cells.Item(2,2).Value = "n.k.gangwar581993@gmail.com"
$wsQuery.Cells.Item(2,3).Value = "niteesh is entering some data in query box user 1"

$wsQuery.Cells.Item(3,1).Value = "pawan"
$wsQuery.Cells.Item(3,2).Value = "pavanoltraining@gmail.com"
$wsQuery.Cells.Item(3,3).Value = "pawan is entering some data in query box user 2"

$wsQuery.Cells.Item(4,1).Value = "rohit"
$wsQuery.Cells.Item(4,2).Value = "rohit11@gmail.com"
$wsQuery.Cells.Item(4,3).Value = "rohit is entering some data in query box user 3"

$wsQuery.Cells.Item(5,1).Value = "kamal"
$wsQuery.Cells.Item(5,2).Value = "kamal12@gmail.com"
$wsQuery.Cells.Item(5,3).Value = "kamal is entering some data in query box user 4"

# ---------------------------------------------------------------------
# 6) Selections / active sheet, matching the final authored state:
#    "login" is the active tab, with a few stray selections left
#    behind on the other sheets from the editing session.
# ---------------------------------------------------------------------
$wsReg.Range("C6").Select()
$wsAddress.Range("A1").Select()
$wsQuery.Range("B2").Select()

$wsLogin.Activate()
$wsLogin.Range("C14").Select()
